$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / percentage / coin-name / URL updates (never numeric-looking)
$ws.Range("D2").Value = '30.520.84'
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").Value = '1.880.53'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("E8").Value = '  +2.46%  '
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("E10").Value = '  +4.56%  '
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("E12").Value = '  +4.17%  '
$ws.Range("E13").Value = '  +9.31%  '
$ws.Range("D14").Value = '1.877.09'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("E15").Value = '  +2.04%  '
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = '30.495.47'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("E18").Value = '  +2.28%  '
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("D21").Value = '2.123.70'
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E26").Value = '  -0.87%  '
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("E31").Value = '  +4.85%  '
$ws.Range("E32").Value = '  +2.83%  '
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("E34").Value = '  +3.20%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E41").Value = '  +4.93%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E42").Value = '  +0.41%  '
$ws.Range("E43").Value = '  +4.88%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("E49").Value = '  +5.17%  '
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("E51").Value = '  +1.86%  '

# Price values that look like plain decimals need to be forced to text
# (leading apostrophe marks them as text; Style reset avoids leaving a stray
# "Text" number-format style behind on the cell).
$ws.Range("D4").Value = "'0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'247.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4770"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.2904"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.06530"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'21.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07730"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'97.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.7411"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Value = "'5.139"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'274.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Value = "'0.000007590"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.9997"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.9995"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'5.260"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'6.192"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'9.341"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'163.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'18.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'1.948"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'1.370"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'0.09970"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'1.519"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'4.321"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'4.074"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.04796"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.129"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.7020"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.01874"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'2.727"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'6.366"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'1.965"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'71.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.4237"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.8396"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.9991"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'102.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'9.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'7.098"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'35.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'924.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.05655"
$ws.Range("D51").Style = "Normal"
